# Scheduling_Entry_Template_Hemmeter.docx - collapse split merge-field runs
# into single runs (removing stray w:proofErr splits) and update the jury
# trial placeholders per the commit "Updated templates and scheduling model
# for jury trial and trial to court".

$d = $word.ActiveDocument

function Replace-Text($range, [string]$find, [string]$replace) {
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# --- Main document body -----------------------------------------------

Replace-Text $d.Content "Case No. {{ case_number }}" "Case No. {{ case_number }}"

Replace-Text $d.Content "{{ defendant.first_name }} {{ defendant.last_name }}," "{{ defendant.first_name }} {{ defendant.last_name }},"

Replace-Text $d.Content " {% if pretrial_scheduled is true %}" " {% if pretrial_scheduled is true %}"

Replace-Text $d.Content " {{ pretrial_date }} between 3:00 PM and 5:00 PM." " {{ pretrial_date }} between 3:00 PM and 5:00 PM."

Replace-Text $d.Content " {{ final_pretrial_date }} at {{ final_pretrial_time }}" " {{ final_pretrial_date }} at {{ final_pretrial_time }}"

Replace-Text $d.Content "Jury Trial on {{ trial_date }} at 8:15 AM" "Jury Trial on {{ jury_trial_date }} at {{ jury_trial_time }}"

# --- Footer (first-page footer holds the real content) -----------------

$footer = $d.Sections.Item(1).Footers.Item(2).Range

Replace-Text $footer "Copies served by Dep. Clerk _______________________________ on the following date ________________________ to:" "Copies served by Dep. Clerk _______________________________ on the following date ________________________ to:"

Replace-Text $footer "{{ defendant.first_name }} {{ defendant.last_name}}: PS" "{{ defendant.first_name }} {{ defendant.last_name}}: PS"

Replace-Text $footer " Scheduling Entry {{ case_number }}" " Scheduling Entry {{ case_number }}"
